$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P": full stats update (rows 2-7), add column H ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 15
$ws1.Range("F2").Value = 9
$ws1.Range("G2").Value = 37.5
$ws1.Range("H2").Value = 7.8
$ws1.Range("D3").Value = 11
$ws1.Range("F3").Value = 21
$ws1.Range("G3").Value = 65.63
$ws1.Range("H3").Value = 7.8
$ws1.Range("D4").Value = 20
$ws1.Range("F4").Value = 14
$ws1.Range("G4").Value = 41.18
$ws1.Range("H4").Value = 8.4
$ws1.Range("D5").Value = 13
$ws1.Range("F5").Value = 8
$ws1.Range("G5").Value = 38.1
$ws1.Range("H5").Value = 7.9
$ws1.Range("D6").Value = 9
$ws1.Range("F6").Value = 26
$ws1.Range("G6").Value = 74.29
$ws1.Range("H6").Value = 7.8
$ws1.Range("D7").Value = 13
$ws1.Range("F7").Value = 8
$ws1.Range("G7").Value = 38.1
$ws1.Range("H7").Value = 7.4

# --- Sheet "Estadisticos 2P": only Reprobados (E) update (rows 2-7) ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 9
$ws2.Range("E3").Value = 21
$ws2.Range("E4").Value = 14
$ws2.Range("E5").Value = 8
$ws2.Range("E6").Value = 26
$ws2.Range("E7").Value = 8

# --- Sheet "Estadisticos Final": same full stats update as 1P ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 15
$ws3.Range("F2").Value = 9
$ws3.Range("G2").Value = 37.5
$ws3.Range("H2").Value = 7.8
$ws3.Range("D3").Value = 11
$ws3.Range("F3").Value = 21
$ws3.Range("G3").Value = 65.63
$ws3.Range("H3").Value = 7.8
$ws3.Range("D4").Value = 20
$ws3.Range("F4").Value = 14
$ws3.Range("G4").Value = 41.18
$ws3.Range("H4").Value = 8.4
$ws3.Range("D5").Value = 13
$ws3.Range("F5").Value = 8
$ws3.Range("G5").Value = 38.1
$ws3.Range("H5").Value = 7.9
$ws3.Range("D6").Value = 9
$ws3.Range("F6").Value = 26
$ws3.Range("G6").Value = 74.29
$ws3.Range("H6").Value = 7.8
$ws3.Range("D7").Value = 13
$ws3.Range("F7").Value = 8
$ws3.Range("G7").Value = 38.1
$ws3.Range("H7").Value = 7.4

# --- Sheet "Rescatables": replace student roster (22 rows instead of 26) ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows.Item(27).Delete()
$ws4.Rows.Item(26).Delete()
$ws4.Rows.Item(25).Delete()
$ws4.Rows.Item(24).Delete()

$ws4.Range("A2").Value = 20330051920263
$ws4.Range("B2").Value = "CARRERA"
$ws4.Range("C2").Value = "ZAVALETA"
$ws4.Range("D2").Value = "ALFREDO"
$ws4.Range("F2").Value = "3APM"
$ws4.Range("A3").Value = 20330051920352
$ws4.Range("B3").Value = "SANCHEZ"
$ws4.Range("C3").Value = "RODRIGUEZ"
$ws4.Range("D3").Value = "EMILIO"
$ws4.Range("F3").Value = "3APM"
$ws4.Range("A4").Value = 20330051920357
$ws4.Range("B4").Value = "XILCAHUA"
$ws4.Range("C4").Value = "TLAXCALA"
$ws4.Range("D4").Value = "LUIS ANGEL"
$ws4.Range("F4").Value = "3APM"
$ws4.Range("A5").Value = 20330051920326
$ws4.Range("B5").Value = "LUNA"
$ws4.Range("C5").Value = "MORALES"
$ws4.Range("D5").Value = "JESUS ANTONIO"
$ws4.Range("F5").Value = "3ASV"
$ws4.Range("A6").Value = 20330051920396
$ws4.Range("B6").Value = "MUÑOZ"
$ws4.Range("C6").Value = "RODRIGUEZ"
$ws4.Range("D6").Value = "VICTOR HUGO"
$ws4.Range("F6").Value = "3ASV"
$ws4.Range("A7").Value = 20330051920040
$ws4.Range("B7").Value = "AMADOR"
$ws4.Range("C7").Value = "PORRAS"
$ws4.Range("D7").Value = "FRANCISCO ALAN"
$ws4.Range("F7").Value = "3BEM"
$ws4.Range("A8").Value = 20330051920295
$ws4.Range("B8").Value = "HERRERA"
$ws4.Range("C8").Value = "CERON"
$ws4.Range("D8").Value = "YAMILE"
$ws4.Range("F8").Value = "3BLCM"
$ws4.Range("A9").Value = 20330051920298
$ws4.Range("B9").Value = "LEYVA"
$ws4.Range("C9").Value = "VELAZQUEZ"
$ws4.Range("D9").Value = "ELIAN"
$ws4.Range("F9").Value = "3BLCM"
$ws4.Range("A10").Value = 20330051920311
$ws4.Range("B10").Value = "ROJAS"
$ws4.Range("C10").Value = "ROJAS"
$ws4.Range("D10").Value = "DULCE MARIA"
$ws4.Range("F10").Value = "3BLCM"
$ws4.Range("A11").Value = 20330051920313
$ws4.Range("B11").Value = "TORRES"
$ws4.Range("C11").Value = "VAZQUEZ"
$ws4.Range("D11").Value = "JOSELIN GUADALUPE"
$ws4.Range("F11").Value = "3BLCM"
$ws4.Range("A12").Value = 19220030050208
$ws4.Range("B12").Value = "DE LUNA"
$ws4.Range("C12").Value = "CORDOVA"
$ws4.Range("D12").Value = "EUDY"
$ws4.Range("F12").Value = "3APV"
$ws4.Range("A13").Value = 20330051920273
$ws4.Range("B13").Value = "LUNA"
$ws4.Range("C13").Value = "FLORES"
$ws4.Range("D13").Value = "MIRANDA"
$ws4.Range("F13").Value = "3APV"
$ws4.Range("A14").Value = 20330051920276
$ws4.Range("B14").Value = "OJEDA"
$ws4.Range("C14").Value = "JIMENEZ"
$ws4.Range("D14").Value = "KAREN YAZMIN"
$ws4.Range("F14").Value = "3APV"
$ws4.Range("A15").Value = 20330051920121
$ws4.Range("B15").Value = "CUATRA"
$ws4.Range("C15").Value = "ZOPIYACTLE"
$ws4.Range("D15").Value = "MARIA"
$ws4.Range("F15").Value = "3ARHV"
$ws4.Range("A16").Value = 20330051920132
$ws4.Range("B16").Value = "LOPEZ"
$ws4.Range("C16").Value = "NOYOLA"
$ws4.Range("D16").Value = "MARIA JOSE"
$ws4.Range("F16").Value = "3ARHV"
$ws4.Range("A17").Value = 20330051920134
$ws4.Range("B17").Value = "MARTINEZ"
$ws4.Range("C17").Value = "CHIPAHUA"
$ws4.Range("D17").Value = "GERMAN ERNESTO"
$ws4.Range("F17").Value = "3ARHV"
$ws4.Range("A18").Value = 20330051920321
$ws4.Range("B18").Value = "COCOTLE"
$ws4.Range("C18").Value = "TLAXCALA"
$ws4.Range("D18").Value = "SURISADAY"
$ws4.Range("F18").Value = "3ASV"
$ws4.Range("A19").Value = 20330051920322
$ws4.Range("B19").Value = "GARCIA"
$ws4.Range("C19").Value = "FLORES"
$ws4.Range("D19").Value = "MARCOS"
$ws4.Range("F19").Value = "3ASV"
$ws4.Range("A20").Value = 20330051920324
$ws4.Range("B20").Value = "GONZALEZ"
$ws4.Range("C20").Value = "HERNANDEZ"
$ws4.Range("D20").Value = "ARIZBETH"
$ws4.Range("F20").Value = "3ASV"
$ws4.Range("A21").Value = 20330051920292
$ws4.Range("B21").Value = "FLORES"
$ws4.Range("C21").Value = "GAMBOA"
$ws4.Range("D21").Value = "VALERIA ANGELY"
$ws4.Range("F21").Value = "3BLCM"
$ws4.Range("A22").Value = 20330051920299
$ws4.Range("B22").Value = "LOPEZ"
$ws4.Range("C22").Value = "MONTERROSAS"
$ws4.Range("D22").Value = "MARIA MAGDALENA"
$ws4.Range("F22").Value = "3BLCM"
$ws4.Range("A23").Value = 20330051920309
$ws4.Range("B23").Value = "RAMOS"
$ws4.Range("C23").Value = "ZEPEDA"
$ws4.Range("D23").Value = "SAMANTHA"
$ws4.Range("F23").Value = "3BLCM"
